$wb = $excel.ActiveWorkbook

# --- 1. "总计" sheet: insert a new 2022-Q4 summary row, shifting existing rows down ---
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 27
$total.Cells.Item(2, 4).Value = 1.37

for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# Match column-A style (s=2) on the freshly inserted row to the rows below it
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)

# --- 2. New "2022-Q4" worksheet: fund-holding detail, placed right before "2022-Q3" ---
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Extend the index-column style (s=2) down to the new rows
$q4.Range("A2").Copy()
$q4.Range("A2:A28").PasteSpecial(-4122)

# Force text storage for the numeric-looking string columns (matches source data typing)
$q4.Range("B2:B28").NumberFormat = "@"
$q4.Range("D2:F28").NumberFormat = "@"
$q4.Range("G2:G27").NumberFormat = "@"

$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "008273"
$q4.Cells.Item(2, 3).Value = "广发优质生活混合"
$q4.Cells.Item(2, 4).Value = "9.28"
$q4.Cells.Item(2, 5).Value = "89.90"
$q4.Cells.Item(2, 6).Value = "3.00"
$q4.Cells.Item(2, 7).Value = "0.2784"
$q4.Cells.Item(2, 8).Value = 9
$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "012690"
$q4.Cells.Item(3, 3).Value = "广发消费领先混合A"
$q4.Cells.Item(3, 4).Value = "5.82"
$q4.Cells.Item(3, 5).Value = "93.14"
$q4.Cells.Item(3, 6).Value = "4.51"
$q4.Cells.Item(3, 7).Value = "0.2625"
$q4.Cells.Item(3, 8).Value = 8
$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Value = "010081"
$q4.Cells.Item(4, 3).Value = "泰康浩泽混合A"
$q4.Cells.Item(4, 4).Value = "6.73"
$q4.Cells.Item(4, 5).Value = "24.06"
$q4.Cells.Item(4, 6).Value = "2.30"
$q4.Cells.Item(4, 7).Value = "0.1548"
$q4.Cells.Item(4, 8).Value = 2
$q4.Cells.Item(5, 1).Value = 3
$q4.Cells.Item(5, 2).Value = "006007"
$q4.Cells.Item(5, 3).Value = "诺安积极配置混合A"
$q4.Cells.Item(5, 4).Value = "3.33"
$q4.Cells.Item(5, 5).Value = "81.11"
$q4.Cells.Item(5, 6).Value = "3.56"
$q4.Cells.Item(5, 7).Value = "0.1185"
$q4.Cells.Item(5, 8).Value = 7
$q4.Cells.Item(6, 1).Value = 4
$q4.Cells.Item(6, 2).Value = "001140"
$q4.Cells.Item(6, 3).Value = "工银总回报灵活配置混合A"
$q4.Cells.Item(6, 4).Value = "4.41"
$q4.Cells.Item(6, 5).Value = "81.48"
$q4.Cells.Item(6, 6).Value = "2.30"
$q4.Cells.Item(6, 7).Value = "0.1014"
$q4.Cells.Item(6, 8).Value = 9
$q4.Cells.Item(7, 1).Value = 5
$q4.Cells.Item(7, 2).Value = "005526"
$q4.Cells.Item(7, 3).Value = "工银瑞信新生代消费灵活配置混合"
$q4.Cells.Item(7, 4).Value = "1.96"
$q4.Cells.Item(7, 5).Value = "90.63"
$q4.Cells.Item(7, 6).Value = "3.34"
$q4.Cells.Item(7, 7).Value = "0.0655"
$q4.Cells.Item(7, 8).Value = 7
$q4.Cells.Item(8, 1).Value = 6
$q4.Cells.Item(8, 2).Value = "001798"
$q4.Cells.Item(8, 3).Value = "泰康新回报灵活配置混合A"
$q4.Cells.Item(8, 4).Value = "1.03"
$q4.Cells.Item(8, 5).Value = "89.15"
$q4.Cells.Item(8, 6).Value = "6.16"
$q4.Cells.Item(8, 7).Value = "0.0634"
$q4.Cells.Item(8, 8).Value = 2
$q4.Cells.Item(9, 1).Value = 7
$q4.Cells.Item(9, 2).Value = "519678"
$q4.Cells.Item(9, 3).Value = "银河消费驱动混合A"
$q4.Cells.Item(9, 4).Value = "0.86"
$q4.Cells.Item(9, 5).Value = "83.68"
$q4.Cells.Item(9, 6).Value = "5.28"
$q4.Cells.Item(9, 7).Value = "0.0454"
$q4.Cells.Item(9, 8).Value = 6
$q4.Cells.Item(10, 1).Value = 8
$q4.Cells.Item(10, 2).Value = "004076"
$q4.Cells.Item(10, 3).Value = "国联安锐意成长混合"
$q4.Cells.Item(10, 4).Value = "1.02"
$q4.Cells.Item(10, 5).Value = "92.87"
$q4.Cells.Item(10, 6).Value = "4.16"
$q4.Cells.Item(10, 7).Value = "0.0424"
$q4.Cells.Item(10, 8).Value = 4
$q4.Cells.Item(11, 1).Value = 9
$q4.Cells.Item(11, 2).Value = "006568"
$q4.Cells.Item(11, 3).Value = "国联安行业领先混合"
$q4.Cells.Item(11, 4).Value = "1.18"
$q4.Cells.Item(11, 5).Value = "93.20"
$q4.Cells.Item(11, 6).Value = "3.22"
$q4.Cells.Item(11, 7).Value = "0.0380"
$q4.Cells.Item(11, 8).Value = 6
$q4.Cells.Item(12, 1).Value = 10
$q4.Cells.Item(12, 2).Value = "015287"
$q4.Cells.Item(12, 3).Value = "永赢优质生活混合A"
$q4.Cells.Item(12, 4).Value = "1.61"
$q4.Cells.Item(12, 5).Value = "36.42"
$q4.Cells.Item(12, 6).Value = "2.18"
$q4.Cells.Item(12, 7).Value = "0.0351"
$q4.Cells.Item(12, 8).Value = 4
$q4.Cells.Item(13, 1).Value = 11
$q4.Cells.Item(13, 2).Value = "001209"
$q4.Cells.Item(13, 3).Value = "前海开源一带一路主题精选灵活配置混合A"
$q4.Cells.Item(13, 4).Value = "0.69"
$q4.Cells.Item(13, 5).Value = "84.87"
$q4.Cells.Item(13, 6).Value = "4.29"
$q4.Cells.Item(13, 7).Value = "0.0296"
$q4.Cells.Item(13, 8).Value = 4
$q4.Cells.Item(14, 1).Value = 12
$q4.Cells.Item(14, 2).Value = "001007"
$q4.Cells.Item(14, 3).Value = "国联安鑫安灵活配置混合"
$q4.Cells.Item(14, 4).Value = "0.74"
$q4.Cells.Item(14, 5).Value = "93.01"
$q4.Cells.Item(14, 6).Value = "3.88"
$q4.Cells.Item(14, 7).Value = "0.0287"
$q4.Cells.Item(14, 8).Value = 7
$q4.Cells.Item(15, 1).Value = 13
$q4.Cells.Item(15, 2).Value = "012691"
$q4.Cells.Item(15, 3).Value = "广发消费领先混合C"
$q4.Cells.Item(15, 4).Value = "0.50"
$q4.Cells.Item(15, 5).Value = "93.14"
$q4.Cells.Item(15, 6).Value = "4.51"
$q4.Cells.Item(15, 7).Value = "0.0226"
$q4.Cells.Item(15, 8).Value = 8
$q4.Cells.Item(16, 1).Value = 14
$q4.Cells.Item(16, 2).Value = "015288"
$q4.Cells.Item(16, 3).Value = "永赢优质生活混合C"
$q4.Cells.Item(16, 4).Value = "0.95"
$q4.Cells.Item(16, 5).Value = "36.42"
$q4.Cells.Item(16, 6).Value = "2.18"
$q4.Cells.Item(16, 7).Value = "0.0207"
$q4.Cells.Item(16, 8).Value = 4
$q4.Cells.Item(17, 1).Value = 15
$q4.Cells.Item(17, 2).Value = "006008"
$q4.Cells.Item(17, 3).Value = "诺安积极配置混合C"
$q4.Cells.Item(17, 4).Value = "0.52"
$q4.Cells.Item(17, 5).Value = "81.11"
$q4.Cells.Item(17, 6).Value = "3.56"
$q4.Cells.Item(17, 7).Value = "0.0185"
$q4.Cells.Item(17, 8).Value = 7
$q4.Cells.Item(18, 1).Value = 16
$q4.Cells.Item(18, 2).Value = "005329"
$q4.Cells.Item(18, 3).Value = "汇添富民安增益定期开放混合A"
$q4.Cells.Item(18, 4).Value = "1.37"
$q4.Cells.Item(18, 5).Value = "28.23"
$q4.Cells.Item(18, 6).Value = "1.02"
$q4.Cells.Item(18, 7).Value = "0.0140"
$q4.Cells.Item(18, 8).Value = 10
$q4.Cells.Item(19, 1).Value = 17
$q4.Cells.Item(19, 2).Value = "010082"
$q4.Cells.Item(19, 3).Value = "泰康浩泽混合C"
$q4.Cells.Item(19, 4).Value = "0.46"
$q4.Cells.Item(19, 5).Value = "24.06"
$q4.Cells.Item(19, 6).Value = "2.30"
$q4.Cells.Item(19, 7).Value = "0.0106"
$q4.Cells.Item(19, 8).Value = 2
$q4.Cells.Item(20, 1).Value = 18
$q4.Cells.Item(20, 2).Value = "001799"
$q4.Cells.Item(20, 3).Value = "泰康新回报灵活配置混合C"
$q4.Cells.Item(20, 4).Value = "0.15"
$q4.Cells.Item(20, 5).Value = "89.15"
$q4.Cells.Item(20, 6).Value = "6.16"
$q4.Cells.Item(20, 7).Value = "0.0092"
$q4.Cells.Item(20, 8).Value = 2
$q4.Cells.Item(21, 1).Value = 19
$q4.Cells.Item(21, 2).Value = "002080"
$q4.Cells.Item(21, 3).Value = "前海开源一带一路主题精选灵活配置混合C"
$q4.Cells.Item(21, 4).Value = "0.10"
$q4.Cells.Item(21, 5).Value = "84.87"
$q4.Cells.Item(21, 6).Value = "4.29"
$q4.Cells.Item(21, 7).Value = "0.0043"
$q4.Cells.Item(21, 8).Value = 4
$q4.Cells.Item(22, 1).Value = 20
$q4.Cells.Item(22, 2).Value = "005330"
$q4.Cells.Item(22, 3).Value = "汇添富民安增益定期开放混合C"
$q4.Cells.Item(22, 4).Value = "0.31"
$q4.Cells.Item(22, 5).Value = "28.23"
$q4.Cells.Item(22, 6).Value = "1.02"
$q4.Cells.Item(22, 7).Value = "0.0032"
$q4.Cells.Item(22, 8).Value = 10
$q4.Cells.Item(23, 1).Value = 21
$q4.Cells.Item(23, 2).Value = "014649"
$q4.Cells.Item(23, 3).Value = "永赢优质精选混合A"
$q4.Cells.Item(23, 4).Value = "0.08"
$q4.Cells.Item(23, 5).Value = "74.97"
$q4.Cells.Item(23, 6).Value = "3.00"
$q4.Cells.Item(23, 7).Value = "0.0024"
$q4.Cells.Item(23, 8).Value = 6
$q4.Cells.Item(24, 1).Value = 22
$q4.Cells.Item(24, 2).Value = "001530"
$q4.Cells.Item(24, 3).Value = "万家瑞富灵活配置混合A"
$q4.Cells.Item(24, 4).Value = "0.21"
$q4.Cells.Item(24, 5).Value = "23.67"
$q4.Cells.Item(24, 6).Value = "0.74"
$q4.Cells.Item(24, 7).Value = "0.0016"
$q4.Cells.Item(24, 8).Value = 6
$q4.Cells.Item(25, 1).Value = 23
$q4.Cells.Item(25, 2).Value = "012007"
$q4.Cells.Item(25, 3).Value = "万家瑞富灵活配置混合C"
$q4.Cells.Item(25, 4).Value = "0.11"
$q4.Cells.Item(25, 5).Value = "23.67"
$q4.Cells.Item(25, 6).Value = "0.74"
$q4.Cells.Item(25, 7).Value = "0.0008"
$q4.Cells.Item(25, 8).Value = 6
$q4.Cells.Item(26, 1).Value = 24
$q4.Cells.Item(26, 2).Value = "015668"
$q4.Cells.Item(26, 3).Value = "银河消费驱动混合C"
$q4.Cells.Item(26, 4).Value = "0.01"
$q4.Cells.Item(26, 5).Value = "83.68"
$q4.Cells.Item(26, 6).Value = "5.28"
$q4.Cells.Item(26, 7).Value = "0.0005"
$q4.Cells.Item(26, 8).Value = 6
$q4.Cells.Item(27, 1).Value = 25
$q4.Cells.Item(27, 2).Value = "011477"
$q4.Cells.Item(27, 3).Value = "工银总回报灵活配置混合C"
$q4.Cells.Item(27, 4).Value = "0.02"
$q4.Cells.Item(27, 5).Value = "81.48"
$q4.Cells.Item(27, 6).Value = "2.30"
$q4.Cells.Item(27, 7).Value = "0.0005"
$q4.Cells.Item(27, 8).Value = 9
$q4.Cells.Item(28, 1).Value = 26
$q4.Cells.Item(28, 2).Value = "014650"
$q4.Cells.Item(28, 3).Value = "永赢优质精选混合C"
$q4.Cells.Item(28, 4).Value = "0.00"
$q4.Cells.Item(28, 5).Value = "74.97"
$q4.Cells.Item(28, 6).Value = "3.00"
$q4.Cells.Item(28, 7).Value = 0
$q4.Cells.Item(28, 8).Value = 6

# Restore the originally active sheet/selection
$total.Activate()
$null = $total.Range("A1").Select()
